$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column J
$ws.Range("J1").Value2 = 'Misc. Data'

$ws.Range("E2").Value2 = '[Martin%Chalumeau%NULL%1,                        Philippe%Bidet%NULL%1,                        Gérard%Lina%NULL%1,                        Mostapha%Mokhtari%NULL%1,                        Marie-Claude%André%NULL%1,                        Dominique%Gendrel%NULL%1,                        Edouard%Bingen%NULL%1,                        Josette%Raymond%j.raymond@svp.ap-hop-paris.fr%1]'
$ws.Range("J2").Value2 = ''

$ws.Range("E3").Value2 = '[Hae-Sung%Nam%NULL%0,                        Mi-Yeon%Yeon%NULL%0,                        Mi-Yeon%Yeon%NULL%0,                        Jung Wan%Park%NULL%0,                        Jung Wan%Park%NULL%0,                        Jee-Young%Hong%NULL%2,                        Jee-Young%Hong%NULL%0,                        Ji Woong%Son%NULL%2,                        Ji Woong%Son%NULL%0]'
$ws.Range("J3").Value2 = ''

$ws.Range("E4").Value2 = '[Michael D.%Christian%NULL%1,                        Mona%Loutfy%NULL%1,                        L. Clifford%McDonald%NULL%2,                        Kenneth F.%Martinez%NULL%1,                        Mariana%Ofner%NULL%1,                        Tom%Wong%NULL%1,                        Tamara%Wallington%NULL%1,                        Wayne L.%Gold%NULL%1,                        Barbara%Mederski%NULL%1,                        Karen%Green%NULL%3,                        Donald E.%Low%NULL%1,                        NULL%NULL%NULL%0]'
$ws.Range("J4").Value2 = ''

$ws.Range("E5").Value2 = '[W. Y.%Kim%xref no email%1,     W.%Choi%xref no email%1,     S.-W.%Park%xref no email%1,     E. B.%Wang%xref no email%1,     W.-J.%Lee%xref no email%1,     Y.%Jee%xref no email%1,     K. S.%Lim%xref no email%1,     H.-J.%Lee%xref no email%1,     S.-M.%Kim%xref no email%1,     S.-O.%Lee%xref no email%1,     S.-H.%Choi%xref no email%1,     Y. S.%Kim%xref no email%1,     J. H.%Woo%xref no email%1,     S.-H.%Kim%xref no email%1]'
$ws.Range("I5").Value2 = ''
$ws.Range("J5").Value2 = 'PUBLISHER: Oxford University Press (OUP)'

$ws.Range("E6").Value2 = '[ J.%Knapp%null%1,                M.A.%Weigand%null%1,                E.%Popp%null%1]'
$ws.Range("J6").Value2 = ''

$ws.Range("E7").Value2 = '[Wei%Liu%NULL%0,                        Fang%Tang%NULL%1,                        Li‐Qun%Fang%NULL%1,                        Sake J.%De Vlas%NULL%1,                        Huai‐Jian%Ma%NULL%1,                        Jie‐Ping%Zhou%NULL%1,                        Caspar W. N.%Looman%NULL%1,                        Jan Hendrik%Richardus%NULL%1,                        Wu‐Chun%Cao%NULL%1]'
$ws.Range("J7").Value2 = ''

$ws.Range("E8").Value2 = '[Mark%Loeb%NULL%1,                        Allison%McGeer%NULL%2,                        Bonnie%Henry%NULL%2,                        Marianna%Ofner%NULL%2,                        David%Rose%NULL%1,                        Tammy%Hlywka%NULL%1,                        Joanne%Levie%NULL%1,                        Jane%McQueen%NULL%1,                        Stephanie%Smith%NULL%1,                        Lorraine%Moss%NULL%1,                        Andrew%Smith%NULL%1,                        Karen%Green%NULL%0,                        Stephen D.%Walter%NULL%1]'
$ws.Range("J8").Value2 = ''

$ws.Range("E9").Value2 = '[Janet%Raboud%NULL%1,                        Altynay%Shigayeva%NULL%1,                        Allison%McGeer%NULL%0,                        Erika%Bontovics%NULL%1,                        Martin%Chapman%NULL%1,                        Denise%Gravel%NULL%1,                        Bonnie%Henry%NULL%0,                        Stephen%Lapinsky%NULL%1,                        Mark%Loeb%NULL%1,                        L. Clifford%McDonald%NULL%0,                        Marianna%Ofner%NULL%0,                        Shirley%Paton%NULL%1,                        Donna%Reynolds%NULL%1,                        Damon%Scales%NULL%1,                        Sandy%Shen%NULL%1,                        Andrew%Simor%NULL%1,                        Thomas%Stewart%NULL%1,                        Mary%Vearncombe%NULL%1,                        Dick%Zoutman%NULL%1,                        Karen%Green%NULL%0,                        Joel Mark%Montgomery%NULL%4,                        Joel Mark%Montgomery%NULL%0]'
$ws.Range("J9").Value2 = ''

$ws.Range("E10").Value2 = '[Jan%Schumacher%xref no email%1,     Stuart A.%Gray%xref no email%1,     Sophie%Michel%xref no email%1,     Roger%Alcock%xref no email%1,     Andrea%Brinker%xref no email%1]'
$ws.Range("I10").Value2 = ''
$ws.Range("J10").Value2 = 'PUBLISHER: Cambridge University Press (CUP)'

$ws.Range("E11").Value2 = '[Hyungoo%Shin%NULL%1,                        Jaehoon%Oh%NULL%1,                        Tae Ho%Lim%NULL%1,                        Hyunggoo%Kang%NULL%1,                        Yeongtak%Song%NULL%1,                        Sanghyun%Lee%NULL%1,                        Abdelouahab%Bellou.%NULL%2,                        Abdelouahab%Bellou.%NULL%0]'
$ws.Range("J11").Value2 = ''

$ws.Range("J12").Value2 = ''
